$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current values for columns D,H,I,J,K,L,M,N,O,P,Q across rows 2-37
# before writing, since this edit permutes row data among rows.
# Use .Value2 (not .Value) -- reading back raw date serials/numbers.
$snapshot = @{}
$snapshot["D"] = @{}
$snapshot["H"] = @{}
$snapshot["I"] = @{}
$snapshot["J"] = @{}
$snapshot["K"] = @{}
$snapshot["L"] = @{}
$snapshot["M"] = @{}
$snapshot["N"] = @{}
$snapshot["O"] = @{}
$snapshot["P"] = @{}
$snapshot["Q"] = @{}
for ($r = 2; $r -le 37; $r++) {
    $snapshot["D"][$r] = $ws.Range("D" + $r).Value2
    $snapshot["H"][$r] = $ws.Range("H" + $r).Value2
    $snapshot["I"][$r] = $ws.Range("I" + $r).Value2
    $snapshot["J"][$r] = $ws.Range("J" + $r).Value2
    $snapshot["K"][$r] = $ws.Range("K" + $r).Value2
    $snapshot["L"][$r] = $ws.Range("L" + $r).Value2
    $snapshot["M"][$r] = $ws.Range("M" + $r).Value2
    $snapshot["N"][$r] = $ws.Range("N" + $r).Value2
    $snapshot["O"][$r] = $ws.Range("O" + $r).Value2
    $snapshot["P"][$r] = $ws.Range("P" + $r).Value2
    $snapshot["Q"][$r] = $ws.Range("Q" + $r).Value2
}

# Row permutation: target row -> source row (captured from snapshot)
$mapping = @{}
$mapping[2] = 4
$mapping[3] = 32
$mapping[4] = 15
$mapping[5] = 10
$mapping[6] = 33
$mapping[7] = 20
$mapping[8] = 19
$mapping[9] = 13
$mapping[10] = 16
$mapping[11] = 21
$mapping[12] = 12
$mapping[13] = 8
$mapping[14] = 18
$mapping[15] = 28
$mapping[16] = 29
$mapping[17] = 17
$mapping[18] = 37
$mapping[19] = 11
$mapping[20] = 23
$mapping[21] = 34
$mapping[22] = 35
$mapping[23] = 31
$mapping[24] = 7
$mapping[25] = 36
$mapping[26] = 25
$mapping[27] = 6
$mapping[28] = 9
$mapping[29] = 3
$mapping[30] = 2
$mapping[31] = 14
$mapping[32] = 5
$mapping[33] = 27
$mapping[34] = 26
$mapping[35] = 22
$mapping[36] = 24
$mapping[37] = 30

foreach ($target in $mapping.Keys) {
    $source = $mapping[$target]
    $ws.Range("D" + $target).Value2 = $snapshot["D"][$source]
    $ws.Range("H" + $target).Value2 = $snapshot["H"][$source]
    $ws.Range("I" + $target).Value2 = $snapshot["I"][$source]
    $ws.Range("J" + $target).Value2 = $snapshot["J"][$source]
    $ws.Range("K" + $target).Value2 = $snapshot["K"][$source]
    $ws.Range("L" + $target).Value2 = $snapshot["L"][$source]
    $ws.Range("M" + $target).Value2 = $snapshot["M"][$source]
    $ws.Range("N" + $target).Value2 = $snapshot["N"][$source]
    $ws.Range("O" + $target).Value2 = $snapshot["O"][$source]
    $ws.Range("P" + $target).Value2 = $snapshot["P"][$source]
    $ws.Range("Q" + $target).Value2 = $snapshot["Q"][$source]
}
